# BOT; UPDATE DATA
# Adds one new day (2020-05-01, serial 43952) of COVID-19 PCR data to the
# "all", "kobe" and "other" sheets, and corrects the running totals on the
# preceding rows of "all" / "kobe" that shift because of it.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet "kobe": correct existing rows 72-78, then insert a new row 79 for
# the new day (pushing the trailing footer row down to 80).
# -----------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("F72").Value = 124
$wsKobe.Range("G72").Value = 115

$wsKobe.Range("F73").Value = 123
$wsKobe.Range("G73").Value = 114

$wsKobe.Range("F74").Value = 127
$wsKobe.Range("G74").Value = 118

$wsKobe.Range("F75").Value = 126
$wsKobe.Range("G75").Value = 117

$wsKobe.Range("F76").Value = 124
$wsKobe.Range("G76").Value = 116

$wsKobe.Range("F77").Value = 129
$wsKobe.Range("G77").Value = 120

$wsKobe.Range("D78").Value = 1
$wsKobe.Range("E78").Value = 258
$wsKobe.Range("F78").Value = 128
$wsKobe.Range("G78").Value = 116

$wsKobe.Rows.Item(79).Insert()
$wsKobe.Range("A79").Value = 43952
$wsKobe.Range("B79").Value = 0
$wsKobe.Range("C79").Value = 1950
$wsKobe.Range("D79").Value = 1
$wsKobe.Range("E79").Value = 259
$wsKobe.Range("F79").Value = 129
$wsKobe.Range("G79").Value = 117
$wsKobe.Range("H79").Value = 12
$wsKobe.Range("I79").Value = 4
$wsKobe.Range("J79").Value = 109

$wsKobe.Activate()
$wsKobe.Range("B80").Select()

# -----------------------------------------------------------------------
# Sheet "other": insert a new row 54 for the new day (pushing the footer
# rows down to 55/56). No other cells change on this sheet.
# -----------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Rows.Item(54).Insert()
$wsOther.Range("A54").Value = 43952
$wsOther.Range("B54").Value = 0
$wsOther.Range("C54").Value = 12
$wsOther.Range("D54").Value = 5
$wsOther.Range("E54").Value = 4
$wsOther.Range("F54").Value = 1
$wsOther.Range("G54").Value = 0
$wsOther.Range("H54").Value = 7

$wsOther.Activate()
$wsOther.Range("F60").Select()

# -----------------------------------------------------------------------
# Sheet "all": correct existing rows 17-23, then insert a new row 24 for
# the new day (pushing the two footnote rows down to 25/26). Done last so
# the workbook's active sheet/tab ends back on "all", as in the original
# file.
# -----------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

$wsAll.Range("D17").Value = 130
$wsAll.Range("E17").Value = 120

$wsAll.Range("C18").Value = 219
$wsAll.Range("D18").Value = 129
$wsAll.Range("E18").Value = 119

$wsAll.Range("C19").Value = 223
$wsAll.Range("D19").Value = 133
$wsAll.Range("E19").Value = 123

$wsAll.Range("C20").Value = 224
$wsAll.Range("D20").Value = 132
$wsAll.Range("E20").Value = 122

$wsAll.Range("C21").Value = 241
$wsAll.Range("D21").Value = 129
$wsAll.Range("E21").Value = 120

$wsAll.Range("C22").Value = 248
$wsAll.Range("D22").Value = 134
$wsAll.Range("E22").Value = 124

$wsAll.Range("B23").Value = 258
$wsAll.Range("C23").Value = 253
$wsAll.Range("D23").Value = 133
$wsAll.Range("E23").Value = 120

$wsAll.Rows.Item(24).Insert()
$wsAll.Range("A24").Value = 43952
$wsAll.Range("B24").Value = 259
$wsAll.Range("C24").Value = 254
$wsAll.Range("D24").Value = 134
$wsAll.Range("E24").Value = 121
$wsAll.Range("F24").Value = 13
$wsAll.Range("G24").Value = 4
$wsAll.Range("H24").Value = 116

$wsAll.Activate()
$wsAll.Range("G27").Select()
